$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.112.19"
$ws.Range("E2").Value = "  +3.41%  "

$ws.Range("D3").Value = "1.693.26"
$ws.Range("E3").Value = "  +0.45%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.05"
$ws.Range("E5").Value = "  +0.47%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.519"
$ws.Range("E6").Value = "  +0.06%  "

$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "24.08"
$ws.Range("E8").Value = "  +5.10%  "

$ws.Range("E9").Value = "  +2.00%  "

$ws.Range("E10").Value = "  +0.34%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0885"
$ws.Range("E11").Value = "  -0.67%  "

$ws.Range("D12").Value = "1.933.59"
$ws.Range("E12").Value = "  +0.50%  "

$ws.Range("D13").Value = "1.694.18"
$ws.Range("E13").Value = "  +0.47%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.18"
$ws.Range("E14").Value = "  -0.26%  "

$ws.Range("E15").Value = "  -0.16%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.96"
$ws.Range("E16").Value = "  +0.14%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "250.96"
$ws.Range("E17").Value = "  +6.51%  "

$ws.Range("D18").Value = "28.076.24"
$ws.Range("E18").Value = "  +3.24%  "

$ws.Range("D19").Value = "0.0₃0743"
$ws.Range("E19").Value = "  +0.11%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.72"
$ws.Range("E20").Value = "  -3.40%  "

$ws.Range("E21").Value = "  -0.07%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.54"
$ws.Range("E22").Value = "  -0.47%  "

$ws.Range("E23").Value = "  -0.31%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.04"
$ws.Range("E24").Value = "  -1.72%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.64"
$ws.Range("E25").Value = "  +0.52%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.35"
$ws.Range("E26").Value = "  -0.13%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.52"
$ws.Range("E27").Value = "  +0.51%  "

$ws.Range("E28").Value = "  +0.25%  "

$ws.Range("E29").Value = "  +0.09%  "

$ws.Range("E30").Value = "  +0.29%  "

$ws.Range("E31").Value = "  +3.22%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.38"
$ws.Range("E32").Value = "  +0.26%  "

$ws.Range("D33").Value = "1.464.61"
$ws.Range("E33").Value = "  -5.04%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.19"
$ws.Range("E34").Value = "  -1.84%  "

$ws.Range("E35").Value = "  -2.91%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.952"
$ws.Range("E36").Value = "  +0.65%  "

$ws.Range("E37").Value = "  +0.43%  "

$ws.Range("E38").Value = "  -2.12%  "

$ws.Range("E39").Value = "  -0.47%  "

$ws.Range("E40").Value = "  -1.79%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "69.36"
$ws.Range("E41").Value = "  +0.33%  "

$ws.Range("E42").Value = "  -0.06%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.60"
$ws.Range("E43").Value = "  -3.02%  "

$ws.Range("D44").Value = "1.838.50"
$ws.Range("E44").Value = "  +0.36%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.24"
$ws.Range("E45").Value = "  -1.11%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.798"
$ws.Range("E46").Value = "  +1.01%  "

$ws.Range("E47").Value = "  +7.17%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "89.49"
$ws.Range("E48").Value = "  -0.73%  "

$ws.Range("D49").Value = "0.0₆0111"
$ws.Range("E49").Value = "  -1.11%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.103"
$ws.Range("E50").Value = "  -0.91%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.03"
$ws.Range("E51").Value = "  -2.81%  "
